# Error Calculations and Plots
# Two data rows were removed from the missing-data worksheet ("RM 232" and
# "SC 92"), and a handful of previously-missing cells in the remaining rows
# were imputed (new numeric values filled in) while a different set of
# cells were newly marked as missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were dropped from the dataset -----------
# Delete from the bottom up so row numbers above the deletion point don't
# shift before we get to them.
$ws.Rows(28).Delete()   # "SC 92"
$ws.Rows(26).Delete()   # "RM 232"

# --- Apply the individual cell imputations / removals ------------------
$ws.Range("C3").Value = 11.2
$ws.Range("D4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("F6").Value = 16.43
$ws.Range("D9").Value = -14.5
$ws.Range("D10").Value = -14.7
$ws.Range("F12").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("F14").Value = 17.76
$ws.Range("F19").Value = 17.81
$ws.Range("F20").ClearContents()
$ws.Range("C21").Value = 12.7
$ws.Range("C23").ClearContents()
$ws.Range("F25").ClearContents()
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 17.44
$ws.Range("C32").Value = 10.5
